$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "NA" page-number note recorded for 2025-04-02 no longer applies;
# clear it back to an empty text cell (leading apostrophe forces a literal
# empty-text value instead of Excel clearing the cell outright).
$ws.Range("C39").Value = "'"
$ws.Range("C39").Style = "Normal"

# Append the newest scraped row for 2025-04-03 (leading apostrophe keeps
# the date stored as literal text, matching the rest of the column,
# instead of Excel re-interpreting it as a date serial number).
$ws.Range("A40").Value = "'2025-04-03"
$ws.Range("A40").Style = "Normal"

$ws.Range("B40").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C40").Value = "NA"
$ws.Range("D40").Value = 1
